# Update the model comparison metrics on the active sheet.
# Values must be written as literal text (matching the original
# inline-string cell type), so a leading apostrophe is used to stop
# Excel's automatic number/currency parsing from converting them.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "'€27,144.34"
$ws.Range("C2").Value = "'€27,459.02"

$ws.Range("B3").Value = "'-1650788376.1205"
$ws.Range("C3").Value = "'-1689285588.1029"

$ws.Range("B4").Value = "'€22,130.19"
$ws.Range("C4").Value = "'€22,360.75"
